$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Relabel the "Statistic (95% CI) ..." column headers to "Percentage (95% CI) ..."
# in both the Initial/Secondary/Post-DETECT assessment tables (header row appears
# twice in this document). wdReplaceAll (2) updates every match in one pass.
$find.Execute("Statistic (95% CI)", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Percentage (95% CI)", 2)
